# Add a new slide ("QUESTIONS?") as the last slide of the deck, using the
# "Title and Content" layout (same layout family used by the other
# text-only slides in this deck, e.g. slideLayout2.xml).
$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# The target slide only has a single content placeholder (idx="1") - no
# title placeholder - so drop the auto-created title shape.
$s.Shapes.Item(1).Delete()

$shp = $s.Shapes.Item(1)
$shp.Name = "Marcador de Posição de Conteúdo 2"

# Position/size (EMU -> points, 1 pt = 12700 EMU)
$shp.Left = 457200 / 12700
$shp.Top = 2348880 / 12700
$shp.Width = 8229600 / 12700
$shp.Height = 1639084 / 12700

$tr = $shp.TextFrame.TextRange
$tr.Text = "QUESTIONS?"
$tr.LanguageID = "pt-PT"
$tr.Font.Size = 96
$tr.ParagraphFormat.Alignment = 2
